$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 27.08707916460632, 24.64519882202148, 28.97104764375426, 24.31772915048765, 254316690)
    ,@(3, 20.61402843980169, 24.11130523681641, 24.39563697346896, 19.36703216355954, 254316690)
    ,@(4, 24.45590385028711, 24.73339462280273, 25.58627035312975, 23.03172322203249, 254316690)
    ,@(5, 27.69358316122026, 28.08698081970215, 28.35334469694572, 26.04622870770109, 254316690)
    ,@(6, 27.62548664053264, 28.90513038635254, 29.61695917782575, 26.90543001169112, 254316690)
    ,@(7, 32.32868316799782, 32.86990356445312, 34.48943871424537, 31.0272722886515, 254316690)
    ,@(8, 35.27128864848864, 34.53318786621094, 36.12964226751073, 34.05632404150255, 254316690)
    ,@(9, 34.40878521960051, 36.938232421875, 38.14075743517561, 34.04387924490865, 254316690)
    ,@(10, 37.58509944808483, 39.25204849243164, 39.78281726711452, 36.48209618137788, 254316690)
    ,@(11, 33.63749938181371, 34.58293151855469, 35.84350875541964, 33.18966210747399, 254316690)
    ,@(12, 33.6284062767466, 37.56754302978516, 37.75908631470836, 33.49515863606636, 254316690)
    ,@(13, 37.30220989214942, 35.22010803222656, 37.39419074396421, 34.56788481360144, 254316690)
    ,@(14, 37.47623864206496, 38.0469856262207, 38.99543008952996, 35.88150385660576, 254316690)
    ,@(15, 40.03819326764636, 40.80524063110352, 42.49105929500441, 39.02670142621704, 254316690)
    ,@(16, 41.6994347438015, 37.60395431518555, 41.91097814770129, 35.95391834708406, 254316690)
    ,@(17, 37.31972701042721, 39.69948577880859, 40.76187837560072, 36.01935950944642, 254316690)
    ,@(18, 44.0129405814935, 43.94467544555664, 44.78090383041818, 42.82685825754679, 254316690)
    ,@(19, 47.1289517292708, 48.24311065673828, 50.26573934288837, 46.51187768417542, 254316690)
    ,@(20, 46.78366020214141, 45.81940460205078, 47.21413186271162, 44.18361097818868, 254316690)
    ,@(21, 49.65815456507334, 51.48197937011719, 54.36033540147756, 49.17410522856788, 254316690)
    ,@(22, 41.92561566267757, 45.6486701965332, 47.74017763759019, 41.37019178344922, 254316690)
    ,@(23, 51.11822285848606, 52.54761505126953, 53.21001509494332, 49.10486971387965, 254316690)
    ,@(24, 49.54175206591385, 42.82304000854492, 49.86543906931052, 42.60433194245146, 254316690)
    ,@(25, 51.02996927484659, 50.33574295043945, 52.85780516408641, 49.04395850310341, 254316690)
    ,@(26, 53.24266514418533, 51.08895111083984, 55.64352911829603, 50.74471027837064, 254316690)
    ,@(27, 56.41962542297789, 57.43899154663086, 57.82900477344377, 55.46230731923109, 254316690)
    ,@(28, 55.1965642243056, 55.48135757446289, 56.02424595913914, 53.47000411460147, 254316690)
    ,@(29, 53.73399416397633, 54.4046630859375, 54.45831782772586, 50.47901487518325, 254316690)
    ,@(30, 47.84274340870019, 46.24289321899414, 48.76849895506412, 44.99357169524637, 254316690)
    ,@(31, 44.28402714625329, 47.33872222900391, 47.83578715031835, 43.65139749751025, 254316690)
    ,@(32, 33.33449602243562, 35.93932342529297, 37.15976691214414, 31.62223472246419, 254316690)
    ,@(33, 37.29452073734611, 40.98632049560547, 41.95977778042102, 36.59656666189482, 254316690)
    ,@(34, 48.05835068213502, 47.45716094970703, 49.02025848792423, 46.58774278636714, 254316690)
    ,@(35, 50.95528754852929, 52.69504165649414, 52.97414604322037, 48.83408995060414, 254316690)
    ,@(36, 44.55481537514859, 43.34708786010742, 46.09021923702836, 41.51209317500254, 254316690)
    ,@(37, 52.2062410781632, 55.23009872436523, 56.16268748023748, 49.76643635649425, 254316690)
    ,@(38, 49.85302816979944, 46.08892440795898, 50.07109917129904, 44.8373842847456, 254316690)
    ,@(39, 38.9072398281281, 41.62347793579102, 42.14950920572883, 38.9072398281281, 254316690)
    ,@(40, 43.39632669119851, 40.27178955078125, 45.82652142735589, 40.00176913968116, 254316690)
    ,@(41, 39.8803388479025, 41.5922737121582, 42.62332771084327, 37.70150894650918, 254316690)
    ,@(42, 35.65959696181734, 36.67170333862305, 37.22197676837867, 31.84699649860791, 254316690)
    ,@(43, 39.87057849531841, 40.4458236694336, 41.89386245276646, 37.78778485084987, 254316690)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $ws.Cells.Item($r, 9).Value = "OTEX CN"
}
